$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-11 is updated from 45221 (2023-10-22)
# to 45224 (2023-10-25) for all rows.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
